$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.804.14"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").Value = "1.872.15"
$ws.Range("E3").Value = "  -0.27%  "

$ws.Range("D4").Value = "'0.9996"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'0.7327"
$ws.Range("E5").Value = "  -1.48%  "

$ws.Range("D6").Value = "'241.31"
$ws.Range("E6").Value = "  -0.43%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").Value = "'0.3138"
$ws.Range("E8").Value = "  -0.79%  "

$ws.Range("D9").Value = "'0.07137"
$ws.Range("E9").Value = "  -0.42%  "

$ws.Range("D10").Value = "'24.43"
$ws.Range("E10").Value = "  -1.62%  "

$ws.Range("D11").Value = "'0.08168"
$ws.Range("E11").Value = "  -3.13%  "

$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "'0.7427"
$ws.Range("E12").Value = "  -1.76%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.886.66"
$ws.Range("E13").Value = "  +0.64%  "

$ws.Range("D14").Value = "'5.347"
$ws.Range("E14").Value = "  -1.71%  "

$ws.Range("D15").Value = "'92.40"
$ws.Range("E15").Value = "  -0.48%  "

$ws.Range("D16").Value = "29.821.05"
$ws.Range("E16").Value = "  -0.43%  "

$ws.Range("D17").Value = "'6.012"
$ws.Range("E17").Value = "  -1.45%  "

$ws.Range("D18").Value = "'248.35"
$ws.Range("E18").Value = "  +1.61%  "

$ws.Range("D19").Value = "'13.40"
$ws.Range("E19").Value = "  -1.90%  "

$ws.Range("D20").Value = "'0.000007812"
$ws.Range("E20").Value = "  -0.37%  "

$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.147.01"
$ws.Range("E21").Value = "  +1.52%  "

$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.25%  "

$ws.Range("E23").Value = "  +0.32%  "

$ws.Range("D24").Value = "'7.773"
$ws.Range("E24").Value = "  -3.00%  "

$ws.Range("D25").Value = "'0.1544"
$ws.Range("E25").Value = "  -1.73%  "

$ws.Range("D26").Value = "'9.203"
$ws.Range("E26").Value = "  -1.43%  "

$ws.Range("D27").Value = "'163.89"
$ws.Range("E27").Value = "  -0.57%  "

$ws.Range("D28").Value = "'18.55"
$ws.Range("E28").Value = "  -0.68%  "

$ws.Range("D29").Value = "'2.016"
$ws.Range("E29").Value = "  -1.26%  "

$ws.Range("D30").Value = "'1.446"
$ws.Range("E30").Value = "  -1.89%  "

$ws.Range("D31").Value = "'4.526"
$ws.Range("E31").Value = "  -2.01%  "

$ws.Range("E32").Value = "  -0.70%  "

$ws.Range("D33").Value = "'4.186"
$ws.Range("E33").Value = "  -2.36%  "

$ws.Range("D34").Value = "'0.05302"
$ws.Range("E34").Value = "  -0.88%  "

$ws.Range("D35").Value = "'1.232"
$ws.Range("E35").Value = "  -0.68%  "

$ws.Range("D36").Value = "'0.7414"
$ws.Range("E36").Value = "  -2.18%  "

$ws.Range("E37").Value = "  +0.11%  "

$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("D39").Value = "'0.01933"
$ws.Range("E39").Value = "  -1.30%  "

$ws.Range("D40").Value = "'2.733"
$ws.Range("E40").Value = "  -0.67%  "

$ws.Range("D41").Value = "'0.4463"
$ws.Range("E41").Value = "  -0.69%  "

$ws.Range("D42").Value = "'5.979"
$ws.Range("E42").Value = "  -2.42%  "

$ws.Range("D43").Value = "'0.8665"
$ws.Range("E43").Value = "  +0.22%  "

$ws.Range("D44").Value = "'71.28"
$ws.Range("E44").Value = "  -2.06%  "

$ws.Range("D45").Value = "1.041.98"
$ws.Range("E45").Value = "  -6.11%  "

$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("D47").Value = "'103.91"
$ws.Range("E47").Value = "  +0.60%  "

$ws.Range("D48").Value = "'1.816"
$ws.Range("E48").Value = "  -1.99%  "

$ws.Range("D49").Value = "'7.429"
$ws.Range("E49").Value = "  -3.75%  "

$ws.Range("D50").Value = "'9.494"
$ws.Range("E50").Value = "  -0.82%  "

$ws.Range("D51").Value = "2.031.81"
$ws.Range("E51").Value = "  +0.88%  "
